# Add a "version" column (with a "version list" sheet) to the nano metadata
# template, matching the commit "Add version 1 everywhere".
#
# Summary of the edit:
#  1. Insert a new worksheet "version list" right after "Export as TSV",
#     containing the single value "1" (a sibling of the other "* list"
#     validation-source sheets).
#  2. Insert a new column A ("version") on "Export as TSV", shifting every
#     existing column (donor_id..data_path) one letter to the right
#     (B..AA instead of A..Z).
#  3. Re-home the existing header-row cell comments so they stay attached to
#     the same logical header (they do not auto-shift with the column
#     insert), and add a new comment on the new A1 ("version") header.
#  4. Add a list data-validation on the new column A, restricted to the
#     "version list" sheet, mirroring the pattern used for the other
#     validated columns.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Export as TSV")

$cols = @("A","B","C","D","E","F","G","H","I","J","K","L","M","N","O","P","Q","R","S","T","U","V","W","X","Y","Z")
$colsAfter = @("B","C","D","E","F","G","H","I","J","K","L","M","N","O","P","Q","R","S","T","U","V","W","X","Y","Z","AA")

# --- Step 1: capture + remove the existing header comments (A1..Z1) -------
# Comments stay bound to their original cell reference when columns are
# inserted, so we must reposition them manually.
$commentTexts = @{}
foreach ($col in $cols) {
    $cell = $ws.Range($col + "1")
    if ($cell.Comment -ne $null) {
        $commentTexts[$col] = $cell.Comment.Text()
        $cell.Comment.Delete()
    }
}

# --- Step 2: insert the new "version" column at A -------------------------
$ws.Columns.Item(1).Insert()

# Give the new header cell the same look (bold, centered, wrap) as the rest
# of row 1 by copying the format from the (now shifted) donor_id header.
$ws.Range("B1").Copy()
$ws.Range("A1").PasteSpecial(-4122)
$excel.CutCopyMode = 0
$ws.Range("A1").Value = "version"

# --- Step 3: re-add the comments, shifted one column to the right ---------
for ($i = 0; $i -lt $cols.Count; $i++) {
    $origCol = $cols[$i]
    $newCol = $colsAfter[$i]
    if ($commentTexts.ContainsKey($origCol)) {
        $ws.Range($newCol + "1").AddComment($commentTexts[$origCol]) | Out-Null
    }
}
$ws.Range("A1").AddComment("Current version of metadata schema. Template provides the correct value.") | Out-Null

# --- Step 4: new "version list" sheet, placed right after "Export as TSV" -
$verSheet = $wb.Worksheets.Add($null, $ws)
$verSheet.Name = "version list"
$verCell = $verSheet.Range("A1")
$verCell.NumberFormat = "@"
$verCell.Value = "1"

# --- Step 5: data validation for the new column, matching the other list- -
# backed columns (list source = the new "version list" sheet).
$verValidation = $ws.Range("A2:A1048576").Validation
$verValidation.Add(3, 1, 1, "='version list'!`$A`$1:`$A`$1")
$verValidation.ErrorTitle = "Value must come from list"
$verValidation.ErrorMessage = "Value must be one of: 1."

Write-Output "done"
